$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'89.289.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.98%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.098.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.67%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.23%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'213.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.00%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'623.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.51%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.373"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -7.41%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.816"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +14.75%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.08%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'3.096.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.97%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +9.45%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.21%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -5.99%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -1.07%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'88.962.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.98%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'32.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.65%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.97%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.105.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.07%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +2.52%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0000213"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.47%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.22%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'424.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.19%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'8.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.50%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.15%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +6.17%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'11.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.91%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'82.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.57%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'3.235.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.76%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.14%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +9.36%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +8.28%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -2.43%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'510.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.09%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -10.13%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'6.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.18%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -2.30%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -4.50%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'22.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.16%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +4.67%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.49%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.10%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.03%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -1.56%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -4.99%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.17%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +5.74%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0699"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +15.12%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -2.06%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'161.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -6.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.12%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.709"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.81%  "
$ws.Range("E51").Style = "Normal"
